$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.591135333333333
$ws.Range("H2").Value = 4.773406
$ws.Range("I2").Value = 0.4707829234247397
$ws.Range("J2").Value = 0.4707829234247397
$ws.Range("M2").Value = 3.483060666666667
$ws.Range("N2").Value = 10.449182
$ws.Range("O2").Value = 0.2527672867110271
$ws.Range("P2").Value = 0.2527672867110271
$ws.Range("Q2").Value = 5.542020894876888
$ws.Range("R2").Value = 49.878188053892
$ws.Range("S2").Value = 0.1189985221839567
$ws.Range("T2").Value = 0.1189985221839567
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.591135333333333
$ws.Range("H3").Value = 4.773406
$ws.Range("I3").Value = 0.4707829234247397
$ws.Range("J3").Value = 0.4707829234247397
$ws.Range("O3").Value = 0.3353267952677969
$ws.Range("P3").Value = 0.335326795267797
$ws.Range("Q3").Value = 7.352170172680665
$ws.Range("R3").Value = 66.16953155412598
$ws.Range("S3").Value = 0.1578661289788226
$ws.Range("T3").Value = 0.1578661289788226
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.591135333333333
$ws.Range("H4").Value = 4.773406
$ws.Range("I4").Value = 0.4707829234247397
$ws.Range("J4").Value = 0.4707829234247397
$ws.Range("M4").Value = 2.773309666666667
$ws.Range("N4").Value = 8.319929
$ws.Range("O4").Value = 0.2012603358768551
$ws.Range("P4").Value = 0.2012603358768551
$ws.Range("Q4").Value = 4.412711000908222
$ws.Range("R4").Value = 39.714399008174
$ws.Range("S4").Value = 0.09474992929355086
$ws.Range("T4").Value = 0.09474992929355089
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.591135333333333
$ws.Range("H5").Value = 4.773406
$ws.Range("I5").Value = 0.4707829234247397
$ws.Range("J5").Value = 0.4707829234247397
$ws.Range("M5").Value = 2.902635666666666
$ws.Range("N5").Value = 8.707906999999999
$ws.Range("O5").Value = 0.2106455821443209
$ws.Range("P5").Value = 0.2106455821443209
$ws.Range("Q5").Value = 4.618486169026887
$ws.Range("R5").Value = 41.56637552124199
$ws.Range("S5").Value = 0.09916834296840953
$ws.Range("T5").Value = 0.09916834296840954
$ws.Range("I6").Value = 0.3035973020998604
$ws.Range("J6").Value = 0.3035973020998604
$ws.Range("M6").Value = 3.483060666666667
$ws.Range("N6").Value = 10.449182
$ws.Range("O6").Value = 0.2527672867110271
$ws.Range("P6").Value = 0.2527672867110271
$ws.Range("Q6").Value = 3.573924431298222
$ws.Range("R6").Value = 32.165319881684
$ws.Range("S6").Value = 0.07673946630456971
$ws.Range("T6").Value = 0.07673946630456972
$ws.Range("I7").Value = 0.3035973020998604
$ws.Range("J7").Value = 0.3035973020998604
$ws.Range("O7").Value = 0.3353267952677969
$ws.Range("P7").Value = 0.335326795267797
$ws.Range("S7").Value = 0.1018043103650954
$ws.Range("T7").Value = 0.1018043103650954
$ws.Range("I8").Value = 0.3035973020998604
$ws.Range("J8").Value = 0.3035973020998604
$ws.Range("M8").Value = 2.773309666666667
$ws.Range("N8").Value = 8.319929
$ws.Range("O8").Value = 0.2012603358768551
$ws.Range("P8").Value = 0.2012603358768551
$ws.Range("Q8").Value = 2.845657920377556
$ws.Range("R8").Value = 25.610921283398
$ws.Range("S8").Value = 0.06110209499192495
$ws.Range("T8").Value = 0.06110209499192496
$ws.Range("I9").Value = 0.3035973020998604
$ws.Range("J9").Value = 0.3035973020998604
$ws.Range("M9").Value = 2.902635666666666
$ws.Range("N9").Value = 8.707906999999999
$ws.Range("O9").Value = 0.2106455821443209
$ws.Range("P9").Value = 0.2106455821443209
$ws.Range("Q9").Value = 2.978357690848222
$ws.Range("R9").Value = 26.805219217634
$ws.Range("S9").Value = 0.06395143043827035
$ws.Range("T9").Value = 0.06395143043827035
$ws.Range("G10").Value = 0.730693
$ws.Range("H10").Value = 2.192079
$ws.Range("I10").Value = 0.2161964349979826
$ws.Range("J10").Value = 0.2161964349979826
$ws.Range("M10").Value = 3.483060666666667
$ws.Range("N10").Value = 10.449182
$ws.Range("O10").Value = 0.2527672867110271
$ws.Range("P10").Value = 0.2527672867110271
$ws.Range("Q10").Value = 2.545048047708667
$ws.Range("R10").Value = 22.905432429378
$ws.Range("S10").Value = 0.05464738627103699
$ws.Range("T10").Value = 0.054647386271037
$ws.Range("G11").Value = 0.730693
$ws.Range("H11").Value = 2.192079
$ws.Range("I11").Value = 0.2161964349979826
$ws.Range("J11").Value = 0.2161964349979826
$ws.Range("O11").Value = 0.3353267952677969
$ws.Range("P11").Value = 0.335326795267797
$ws.Range("Q11").Value = 3.376318259951
$ws.Range("R11").Value = 30.386864339559
$ws.Range("S11").Value = 0.07249645769619607
$ws.Range("T11").Value = 0.07249645769619609
$ws.Range("G12").Value = 0.730693
$ws.Range("H12").Value = 2.192079
$ws.Range("I12").Value = 0.2161964349979826
$ws.Range("J12").Value = 0.2161964349979826
$ws.Range("M12").Value = 2.773309666666667
$ws.Range("N12").Value = 8.319929
$ws.Range("O12").Value = 0.2012603358768551
$ws.Range("P12").Value = 0.2012603358768551
$ws.Range("Q12").Value = 2.026437960265667
$ws.Range("R12").Value = 18.237941642391
$ws.Range("S12").Value = 0.04351176712307265
$ws.Range("T12").Value = 0.04351176712307266
$ws.Range("G13").Value = 0.730693
$ws.Range("H13").Value = 2.192079
$ws.Range("I13").Value = 0.2161964349979826
$ws.Range("J13").Value = 0.2161964349979826
$ws.Range("M13").Value = 2.902635666666666
$ws.Range("N13").Value = 8.707906999999999
$ws.Range("O13").Value = 0.2106455821443209
$ws.Range("P13").Value = 0.2106455821443209
$ws.Range("Q13").Value = 2.120935563183667
$ws.Range("R13").Value = 19.088420068653
$ws.Range("S13").Value = 0.04554082390767688
$ws.Range("T13").Value = 0.04554082390767688
$ws.Range("G14").Value = 0.03184866666666667
$ws.Range("H14").Value = 0.09554600000000001
$ws.Range("I14").Value = 0.009423339477417213
$ws.Range("J14").Value = 0.009423339477417213
$ws.Range("M14").Value = 3.483060666666667
$ws.Range("N14").Value = 10.449182
$ws.Range("O14").Value = 0.2527672867110271
$ws.Range("P14").Value = 0.2527672867110271
$ws.Range("Q14").Value = 0.1109308381524445
$ws.Range("R14").Value = 0.9983775433720001
$ws.Range("S14").Value = 0.002381911951463657
$ws.Range("T14").Value = 0.002381911951463657
$ws.Range("G15").Value = 0.03184866666666667
$ws.Range("H15").Value = 0.09554600000000001
$ws.Range("I15").Value = 0.009423339477417213
$ws.Range("J15").Value = 0.009423339477417213
$ws.Range("O15").Value = 0.3353267952677969
$ws.Range("P15").Value = 0.335326795267797
$ws.Range("Q15").Value = 0.1471633570073333
$ws.Range("R15").Value = 1.324470213066
$ws.Range("S15").Value = 0.00315989822768283
$ws.Range("T15").Value = 0.00315989822768283
$ws.Range("G16").Value = 0.03184866666666667
$ws.Range("H16").Value = 0.09554600000000001
$ws.Range("I16").Value = 0.009423339477417213
$ws.Range("J16").Value = 0.009423339477417213
$ws.Range("M16").Value = 2.773309666666667
$ws.Range("N16").Value = 8.319929
$ws.Range("O16").Value = 0.2012603358768551
$ws.Range("P16").Value = 0.2012603358768551
$ws.Range("Q16").Value = 0.08832621513711111
$ws.Range("R16").Value = 0.7949359362340001
$ws.Range("S16").Value = 0.001896544468306616
$ws.Range("T16").Value = 0.001896544468306617
$ws.Range("G17").Value = 0.03184866666666667
$ws.Range("H17").Value = 0.09554600000000001
$ws.Range("I17").Value = 0.009423339477417213
$ws.Range("J17").Value = 0.009423339477417213
$ws.Range("M17").Value = 2.902635666666666
$ws.Range("N17").Value = 8.707906999999999
$ws.Range("O17").Value = 0.2106455821443209
$ws.Range("P17").Value = 0.2106455821443209
$ws.Range("Q17").Value = 0.09244507580244445
$ws.Range("R17").Value = 0.8320056822219999
$ws.Range("S17").Value = 0.001984984829964109
$ws.Range("T17").Value = 0.001984984829964109
